$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.873.74"
$ws.Range("E2").Value = "  +4.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.844.84"
$ws.Range("E3").Value = "  +6.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "423.56"
$ws.Range("E5").Value = "  +4.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.77"
$ws.Range("E6").Value = "  -1.93%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.838.66"
$ws.Range("E7").Value = "  +6.00%  "
$ws.Range("E8").Value = "  -1.62%  "
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("E11").Value = "  -1.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000335"
$ws.Range("E12").Value = "  +5.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.73"
$ws.Range("E13").Value = "  -2.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.29"
$ws.Range("E14").Value = "  +3.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.466.76"
$ws.Range("E15").Value = "  +6.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.80"
$ws.Range("E16").Value = "  +17.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.832.94"
$ws.Range("E17").Value = "  +5.89%  "
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.81"
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.137.75"
$ws.Range("E20").Value = "  +4.36%  "
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "415.72"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.97"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.12"
$ws.Range("E24").Value = "  -1.33%  "
$ws.Range("E25").Value = "  +1.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "37.41"
$ws.Range("E26").Value = "  +5.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.89"
$ws.Range("E27").Value = "  +5.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.24"
$ws.Range("E28").Value = "  +1.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.38"
$ws.Range("E29").Value = "  +35.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.34"
$ws.Range("E30").Value = "  +3.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "738.03"
$ws.Range("E31").Value = "  +8.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.11"
$ws.Range("E32").Value = "  +2.92%  "
$ws.Range("E33").Value = "  +5.17%  "
$ws.Range("E34").Value = "  +5.43%  "
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.152"
$ws.Range("E36").Value = "  -4.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "38.59"
$ws.Range("E37").Value = "  -6.02%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.50"
$ws.Range("E38").Value = "  +24.84%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "55.49"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0729"
$ws.Range("E40").Value = "  +13.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0459"
$ws.Range("E41").Value = "  -1.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.88"
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("E43").Value = "  +0.93%  "
$ws.Range("E44").Value = "  +2.40%  "
$ws.Range("E45").Value = "  -4.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.315"
$ws.Range("E46").Value = "  +8.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.10"
$ws.Range("E47").Value = "  -1.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.06"
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.52"
$ws.Range("E49").Value = "  -2.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.80"
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("E51").Value = "  +0.57%  "
